# Auto-generated edit script: applies the 'remove row 2 / shift rows up' change
# to sheets '展览' (sheet1) and '全部类型' (sheet4).
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '2024-04-13'
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = '南昌·原X穹X崩only'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = '丰和北大道299号 新吉花园酒店'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '2024.04.13 10:00-04.13 17:00'
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 166
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 65
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=80807'
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg'
$ws.Range("I2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '2024-04-13'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = '南昌·第二届漫拥动漫嘉年华mini'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '2024.04.13 10:00-04.14 18:00'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 151
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 55
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=82210'
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png'
$ws.Range("I3").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '2024-04-20'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = '南昌·DSL国风动漫游戏嘉年华'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = '九龙大道1388号（上饶街与九龙大道交叉口西北100米） 中国南昌虚拟现实VR产业基地'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '2024.04.20 09:00-04.21 17:00'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 160
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = 55
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=82107'
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = '//i0.hdslb.com/bfs/openplatform/202402/QDlumVb41708943318282.jpeg'
$ws.Range("I4").Style = "Normal"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '2024-04-20'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = '南昌·New World国潮动漫博览会'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '2024.04.20 09:30-04.21 17:00'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = 4658
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = 60
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=82411'
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202403/xbYbLXc81709707724935.jpeg'
$ws.Range("I5").Style = "Normal"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '2024-04-20'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = '南昌·晨啼漫拥二次元随机舞蹈派对-热爱欢聚(免费活动)'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = '莲塘镇澄湖东路1111号 玺悦城生活广场'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '2024.04.20 15:00-04.20 19:00'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = 20
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = 30.99
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=83272'
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202403/wZUteBVO1710507652186.png'
$ws.Range("I6").Style = "Normal"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '2024-04-20'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '抚州·四月之约动漫游戏聚会'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '迎宾大道688号 抚州万达广场'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '2024.04.20 10:00-04.20 17:00'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = 39
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = 20
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=83316'
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = '//i0.hdslb.com/bfs/openplatform/202403/2A7apu3o1711082007471.jpeg'
$ws.Range("I7").Style = "Normal"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '2024-05-01'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '九江·第三届ACD动漫游戏嘉年华'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '九瑞大道与重庆路交汇处西南角 九江国际会展中心'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '2024.05.01 09:00-05.02 17:00'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = 514
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = 55
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=82464'
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg'
$ws.Range("I8").Style = "Normal"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '2024-05-01'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '吉安·COMIC LIFE次元假日04'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '2024.05.01 09:00-05.01 18:00'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = 472
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = 50
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=82949'
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Value = '//i0.hdslb.com/bfs/openplatform/202403/XDTuzcBV1710478728595.jpeg'
$ws.Range("I9").Style = "Normal"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '2024-05-01'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = '景德镇·第一届国际动漫节'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '2024.05.01 10:00-05.02 18:00'
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = 20
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = 55
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=83734'
$ws.Range("H10").Style = "Normal"
$ws.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202404/6Izq7ZEk1712054058852.jpeg'
$ws.Range("I10").Style = "Normal"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '2024-05-01'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = '景德镇·第一届国际动漫节吕书君内场票'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '2024.05.01 10:00-05.01 18:00'
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = 15
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = 128
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=83666'
$ws.Range("H11").Style = "Normal"
$ws.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202403/OBgwOCB21711786949138.jpeg'
$ws.Range("I11").Style = "Normal"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '2024-05-01'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = '景德镇·第十四届瓷都ACG动漫游戏博览会'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = '新厂西路315号 陶溪川发布大厅'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '2024.05.01 10:00-05.02 18:00'
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = 1328
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = 50
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=83016'
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").Value = '//i2.hdslb.com/bfs/openplatform/202403/c0q8seJL1710835930052.png'
$ws.Range("I12").Style = "Normal"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '2024-05-01'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = '江西·广电·Unlimited Project 动漫游戏博览会'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '2024.05.01 09:00-05.02 17:00'
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = 2715
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = 68
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=83230'
$ws.Range("H13").Style = "Normal"
$ws.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202404/S1nqZf721712025221477.jpeg'
$ws.Range("I13").Style = "Normal"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '2024-05-01'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '江西·第二十二届九江ACJJ国际动漫展'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '体育路九江市体育中心-体育馆 九江市体育中心'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '2024.05.01 09:00-05.02 17:00'
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = 384
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = 50
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=83004'
$ws.Range("H14").Style = "Normal"
$ws.Range("I14").Value = '//i2.hdslb.com/bfs/openplatform/202403/lFThDvkh1710829330909.jpeg'
$ws.Range("I14").Style = "Normal"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '2024-05-01'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '赣州·COMIC WORLD次元创作同人季特典·SP·动漫游戏嘉年华'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '2024.05.01 10:00-05.03 17:00'
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = 85
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = 48
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=83537'
$ws.Range("H15").Style = "Normal"
$ws.Range("I15").Value = '//i2.hdslb.com/bfs/openplatform/202403/BKJfMXXx1711091647172.jpeg'
$ws.Range("I15").Style = "Normal"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '2024-05-01'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '赣州·十万伏特-第六届青年文化综合展览会'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '石楠路与仓背岭路交叉口北120米 新旅中书艺术综合体'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '2024.05.01 09:30-05.03 17:00'
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = 66
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = 55
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=83628'
$ws.Range("H16").Style = "Normal"
$ws.Range("I16").Value = '//i1.hdslb.com/bfs/openplatform/202404/k8EDbMuk1711961223856.jpeg'
$ws.Range("I16").Style = "Normal"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '2024-05-02'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '抚州·临次元07国漫&运动番嘉年华'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '伍塘路1098号 乐课篮球公园'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '2024.05.02 10:00-05.02 16:00'
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = 64
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = 50
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=82911'
$ws.Range("H17").Style = "Normal"
$ws.Range("I17").Value = '//i0.hdslb.com/bfs/openplatform/202403/Pc5TMZ001710429899787.jpeg'
$ws.Range("I17").Style = "Normal"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '2024-05-02'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '江西·ShiningStaR数字互娱嘉年华'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '2024.05.02 09:30-05.04 17:00'
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Value = 2331
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = 65
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=83180'
$ws.Range("H18").Style = "Normal"
$ws.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202403/EqmGU5NC1711015780862.jpeg'
$ws.Range("I18").Style = "Normal"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '2024-05-03'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '上饶·第一届星光次元国风动漫游戏嘉年华暨我和我的cos小伙伴们'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '带湖路66-68号 华都臻悦酒店'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '2024.05.03 09:30-05.03 17:30'
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").Value = 98
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = 45
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=83478'
$ws.Range("H19").Style = "Normal"
$ws.Range("I19").Value = '//i2.hdslb.com/bfs/openplatform/202403/GFCaK00i1711614562687.jpeg'
$ws.Range("I19").Style = "Normal"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '2024-05-03'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '新余·LD02国风动漫嘉年华'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = '劳动北路888号 金联体育篮球馆'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '2024.05.03 10:00-05.03 17:00'
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").Value = 75
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = 50
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=83379'
$ws.Range("H20").Style = "Normal"
$ws.Range("I20").Value = '//i0.hdslb.com/bfs/openplatform/202403/jozduadT1711362183223.jpeg'
$ws.Range("I20").Style = "Normal"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '2024-05-03'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '江西·2024南昌玛雅《次元之芯》主题动漫嘉年华'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = '真君路999号 南昌玛雅乐园'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '2024.05.03 10:30-05.04 19:30'
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Value = 30
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = '不可售'
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=83608'
$ws.Range("H21").Style = "Normal"
$ws.Range("I21").Value = '//i2.hdslb.com/bfs/openplatform/202404/83wvFhen1712040649705.jpeg'
$ws.Range("I21").Style = "Normal"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '2024-05-03'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '2024.05.03 09:30-05.03 17:30'
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").Value = 163
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").Value = 188
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=83497'
$ws.Range("H22").Style = "Normal"
$ws.Range("I22").Value = '//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg'
$ws.Range("I22").Style = "Normal"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '2024-05-04'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '2024.05.04 09:30-05.04 17:30'
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Value = 110
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Value = 228
$ws.Range("G23").Style = "Normal"
$ws.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=83593'
$ws.Range("H23").Style = "Normal"
$ws.Range("I23").Value = '//i0.hdslb.com/bfs/openplatform/202404/LcnCzDxF1711935576170.jpeg'
$ws.Range("I23").Style = "Normal"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '2024-05-18'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '南昌·花绒万兽首届兽聚'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '2024.05.18 09:30-05.19 16:30'
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").Value = 45
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Value = 50
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=83689'
$ws.Range("H24").Style = "Normal"
$ws.Range("I24").Value = '//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg'
$ws.Range("I24").Style = "Normal"

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2024-05-26'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '南昌·代号鸢盛花行only'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '2024.05.26 09:30-05.26 17:30'
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").Value = 228
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value = 78
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=82529'
$ws.Range("H25").Style = "Normal"
$ws.Range("I25").Value = '//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png'
$ws.Range("I25").Style = "Normal"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '2024-06-10'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '上饶·ETI动漫节'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = '滨江东路与体育馆路交叉口西100米 力加体育综合运动中心'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '2024.06.10 10:00-06.10 16:00'
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").Value = 40
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value = 36.6
$ws.Range("G26").Style = "Normal"
$ws.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=83422'
$ws.Range("H26").Style = "Normal"
$ws.Range("I26").Value = '//i1.hdslb.com/bfs/openplatform/202403/vvJKFJal1711460768984.jpeg'
$ws.Range("I26").Style = "Normal"

$ws.Rows(27).Delete()

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '2024-04-13'
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = '南昌·原X穹X崩only'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = '丰和北大道299号 新吉花园酒店'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '2024.04.13 10:00-04.13 17:00'
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 166
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 65
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=80807'
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg'
$ws.Range("I2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '2024-04-13'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = '南昌·第二届漫拥动漫嘉年华mini'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '2024.04.13 10:00-04.14 18:00'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 151
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 55
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=82210'
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png'
$ws.Range("I3").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '2024-04-20'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = '南昌·DSL国风动漫游戏嘉年华'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = '九龙大道1388号（上饶街与九龙大道交叉口西北100米） 中国南昌虚拟现实VR产业基地'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '2024.04.20 09:00-04.21 17:00'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 160
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = 55
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=82107'
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = '//i0.hdslb.com/bfs/openplatform/202402/QDlumVb41708943318282.jpeg'
$ws.Range("I4").Style = "Normal"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '2024-04-20'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = '南昌·Kpop New Life'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = '桃苑东路1号桃苑大厦1楼 星辰LIVE SHOW(南昌店)'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '2024.04.20 14:00-04.20 18:00'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = 40
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = 79
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=83625'
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202403/QUP5gL211711892792092.jpeg'
$ws.Range("I5").Style = "Normal"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '2024-04-20'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = '南昌·New World国潮动漫博览会'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '2024.04.20 09:30-04.21 17:00'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = 4658
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = 60
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=82411'
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202403/xbYbLXc81709707724935.jpeg'
$ws.Range("I6").Style = "Normal"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '2024-04-20'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '南昌·晨啼漫拥二次元随机舞蹈派对-热爱欢聚(免费活动)'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '莲塘镇澄湖东路1111号 玺悦城生活广场'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '2024.04.20 15:00-04.20 19:00'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = 20
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = 30.99
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=83272'
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202403/wZUteBVO1710507652186.png'
$ws.Range("I7").Style = "Normal"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '2024-04-20'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '抚州·四月之约动漫游戏聚会'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '迎宾大道688号 抚州万达广场'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '2024.04.20 10:00-04.20 17:00'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = 39
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = 20
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=83316'
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202403/2A7apu3o1711082007471.jpeg'
$ws.Range("I8").Style = "Normal"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '2024-05-01'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '九江·第三届ACD动漫游戏嘉年华'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '九瑞大道与重庆路交汇处西南角 九江国际会展中心'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '2024.05.01 09:00-05.02 17:00'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = 514
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = 55
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=82464'
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Value = '//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg'
$ws.Range("I9").Style = "Normal"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '2024-05-01'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = '吉安·COMIC LIFE次元假日04'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '2024.05.01 09:00-05.01 18:00'
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = 472
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = 50
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=82949'
$ws.Range("H10").Style = "Normal"
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202403/XDTuzcBV1710478728595.jpeg'
$ws.Range("I10").Style = "Normal"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '2024-05-01'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = '景德镇·第一届国际动漫节'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '2024.05.01 10:00-05.02 18:00'
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = 20
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = 55
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=83734'
$ws.Range("H11").Style = "Normal"
$ws.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202404/6Izq7ZEk1712054058852.jpeg'
$ws.Range("I11").Style = "Normal"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '2024-05-01'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = '景德镇·第一届国际动漫节吕书君内场票'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '2024.05.01 10:00-05.01 18:00'
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = 15
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = 128
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=83666'
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202403/OBgwOCB21711786949138.jpeg'
$ws.Range("I12").Style = "Normal"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '2024-05-01'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = '景德镇·第十四届瓷都ACG动漫游戏博览会'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '新厂西路315号 陶溪川发布大厅'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '2024.05.01 10:00-05.02 18:00'
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = 1328
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = 50
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=83016'
$ws.Range("H13").Style = "Normal"
$ws.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202403/c0q8seJL1710835930052.png'
$ws.Range("I13").Style = "Normal"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '2024-05-01'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '江西·广电·Unlimited Project 动漫游戏博览会'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '2024.05.01 09:00-05.02 17:00'
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = 2715
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = 68
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=83230'
$ws.Range("H14").Style = "Normal"
$ws.Range("I14").Value = '//i2.hdslb.com/bfs/openplatform/202404/S1nqZf721712025221477.jpeg'
$ws.Range("I14").Style = "Normal"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '2024-05-01'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '江西·第二十二届九江ACJJ国际动漫展'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '体育路九江市体育中心-体育馆 九江市体育中心'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '2024.05.01 09:00-05.02 17:00'
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = 384
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = 50
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=83004'
$ws.Range("H15").Style = "Normal"
$ws.Range("I15").Value = '//i2.hdslb.com/bfs/openplatform/202403/lFThDvkh1710829330909.jpeg'
$ws.Range("I15").Style = "Normal"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '2024-05-01'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '赣州·COMIC WORLD次元创作同人季特典·SP·动漫游戏嘉年华'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '2024.05.01 10:00-05.03 17:00'
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = 85
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = 48
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=83537'
$ws.Range("H16").Style = "Normal"
$ws.Range("I16").Value = '//i2.hdslb.com/bfs/openplatform/202403/BKJfMXXx1711091647172.jpeg'
$ws.Range("I16").Style = "Normal"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '2024-05-01'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '赣州·十万伏特-第六届青年文化综合展览会'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '石楠路与仓背岭路交叉口北120米 新旅中书艺术综合体'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '2024.05.01 09:30-05.03 17:00'
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = 66
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = 55
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=83628'
$ws.Range("H17").Style = "Normal"
$ws.Range("I17").Value = '//i1.hdslb.com/bfs/openplatform/202404/k8EDbMuk1711961223856.jpeg'
$ws.Range("I17").Style = "Normal"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '2024-05-02'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '抚州·临次元07国漫&运动番嘉年华'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '伍塘路1098号 乐课篮球公园'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '2024.05.02 10:00-05.02 16:00'
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Value = 64
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = 50
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=82911'
$ws.Range("H18").Style = "Normal"
$ws.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202403/Pc5TMZ001710429899787.jpeg'
$ws.Range("I18").Style = "Normal"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '2024-05-02'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '江西·ShiningStaR数字互娱嘉年华'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '2024.05.02 09:30-05.04 17:00'
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").Value = 2331
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = 65
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=83180'
$ws.Range("H19").Style = "Normal"
$ws.Range("I19").Value = '//i0.hdslb.com/bfs/openplatform/202403/EqmGU5NC1711015780862.jpeg'
$ws.Range("I19").Style = "Normal"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '2024-05-03'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '上饶·第一届星光次元国风动漫游戏嘉年华暨我和我的cos小伙伴们'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = '带湖路66-68号 华都臻悦酒店'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '2024.05.03 09:30-05.03 17:30'
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").Value = 98
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = 45
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=83478'
$ws.Range("H20").Style = "Normal"
$ws.Range("I20").Value = '//i2.hdslb.com/bfs/openplatform/202403/GFCaK00i1711614562687.jpeg'
$ws.Range("I20").Style = "Normal"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '2024-05-03'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '新余·LD02国风动漫嘉年华'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = '劳动北路888号 金联体育篮球馆'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '2024.05.03 10:00-05.03 17:00'
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Value = 75
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = 50
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=83379'
$ws.Range("H21").Style = "Normal"
$ws.Range("I21").Value = '//i0.hdslb.com/bfs/openplatform/202403/jozduadT1711362183223.jpeg'
$ws.Range("I21").Style = "Normal"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '2024-05-03'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '江西·2024南昌玛雅《次元之芯》主题动漫嘉年华'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = '真君路999号 南昌玛雅乐园'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '2024.05.03 10:30-05.04 19:30'
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").Value = 30
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").Value = '不可售'
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=83608'
$ws.Range("H22").Style = "Normal"
$ws.Range("I22").Value = '//i2.hdslb.com/bfs/openplatform/202404/83wvFhen1712040649705.jpeg'
$ws.Range("I22").Style = "Normal"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '2024-05-03'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '2024.05.03 09:30-05.03 17:30'
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Value = 163
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Value = 188
$ws.Range("G23").Style = "Normal"
$ws.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=83497'
$ws.Range("H23").Style = "Normal"
$ws.Range("I23").Value = '//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg'
$ws.Range("I23").Style = "Normal"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '2024-05-04'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '2024.05.04 09:30-05.04 17:30'
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").Value = 110
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Value = 228
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=83593'
$ws.Range("H24").Style = "Normal"
$ws.Range("I24").Value = '//i0.hdslb.com/bfs/openplatform/202404/LcnCzDxF1711935576170.jpeg'
$ws.Range("I24").Style = "Normal"

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2024-05-18'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '南昌·花绒万兽首届兽聚'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '2024.05.18 09:30-05.19 16:30'
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").Value = 45
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value = 50
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=83689'
$ws.Range("H25").Style = "Normal"
$ws.Range("I25").Value = '//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg'
$ws.Range("I25").Style = "Normal"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '2024-05-26'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '南昌·代号鸢盛花行only'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '2024.05.26 09:30-05.26 17:30'
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").Value = 228
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value = 78
$ws.Range("G26").Style = "Normal"
$ws.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=82529'
$ws.Range("H26").Style = "Normal"
$ws.Range("I26").Value = '//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png'
$ws.Range("I26").Style = "Normal"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '2024-06-10'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = '上饶·ETI动漫节'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = '滨江东路与体育馆路交叉口西100米 力加体育综合运动中心'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '2024.06.10 10:00-06.10 16:00'
$ws.Range("E27").Style = "Normal"
$ws.Range("F27").Value = 40
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Value = 36.6
$ws.Range("G27").Style = "Normal"
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=83422'
$ws.Range("H27").Style = "Normal"
$ws.Range("I27").Value = '//i1.hdslb.com/bfs/openplatform/202403/vvJKFJal1711460768984.jpeg'
$ws.Range("I27").Style = "Normal"

$ws.Rows(28).Delete()
